# Regenerate save_data: replace the old "Strike#" counts in column G (K)
# with the recalculated strikeout totals pulled from the boxscore data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 6
    4  = 3
    5  = 0
    6  = 5
    7  = 4
    8  = 4
    9  = 3
    10 = 3
    11 = 5
    12 = 5
    13 = 6
    14 = 2
    15 = 2
    16 = 6
    17 = 3
    18 = 2
    19 = 3
    20 = 4
    21 = 0
    22 = 1
    23 = 2
    24 = 0
    25 = 2
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 3
    31 = 1
    32 = 1
    33 = 2
    34 = 0
    35 = 0
    36 = 1
    37 = 3
    38 = 1
    39 = 2
    40 = 1
    41 = 0
    42 = 1
    44 = 0
    45 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
